# Sample Project / Main.xlsx - "Rules" sheet.
# Cell B11 held the text "R40"; it is retyped as the literal text "1".
# The leading apostrophe is the standard Excel text-prefix so the
# digit-only entry is stored as a string (shared-string "1") instead of
# being auto-converted to the number 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("B11").Value = "'1"
